$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 used to hold the redirect "/contactus" (A4) next to a google.com
# hyperlink (B4, using the built-in "Hyperlink" cell style). The row was
# reworked so A4 now holds the full redirect destination URL and is itself
# turned into a hyperlink (matching B4's existing "Hyperlink" styling),
# while B4 keeps the google URL.
$ws.Range("A4").Value = "https://main--franklindemo--dhanashrideshpande.hlx.live/contactus"

# Turn A4 into a hyperlink pointing at that same URL, then (re)apply the
# "Hyperlink" cell style -- Add() stamps its own formatting on the cell
# first, so Style is (re)applied afterwards to match B4 exactly.
$ws.Hyperlinks.Add($ws.Range("A4"), "https://main--franklindemo--dhanashrideshpande.hlx.live/contactus") | Out-Null
$ws.Range("A4").Style = "Hyperlink"

# Column A is now wide enough to show the full URL, similar to column B's
# existing best-fit width.
$ws.Columns.Item(1).ColumnWidth = 57.498697916666664

# The active selection moved from B14 to A14.
$ws.Range("A14").Select() | Out-Null
